$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-10 (Q0..Q8) with new computed values
$data = @(
    @(6,  -0.02759765538578432, 0.5993211969665078, 0.7764809546060641, 0.881181567332218,  0.8895131873692134, 51),
    @(7,   0.09824497858667457, 0.6346353876862429, 0.9122639460283238, 0.9551250944396361, 0.9597043923121814, 50),
    @(8,  -0.001319298517514361,0.6476906327130449, 0.7146872890618677, 0.8453917961879378, 0.8541515268229536, 49),
    @(9,   0.1228341545534477,  0.6897645129407867, 0.839119603502874,  0.9160347174113402, 0.917367950390179,  48),
    @(10,  0.04626027039139741, 0.6022393271950429, 0.7064898076028564, 0.8405294805078858, 0.8483287965471219, 47),
    @(11,  0.1135060740419775,  0.7098892706801739, 0.8783339292678302, 0.9371947125692879, 0.940575637214629,  46),
    @(12,  0.06286716272461017, 0.7192557759049942, 0.9934777917745472, 0.9967335610756504, 1.005989427279836,  45),
    @(13,  0.1466813223655823,  0.7501106339929442, 0.9202103855659967, 0.959275969450917,  0.9589550867376859, 44),
    @(14,  0.07202425349984963, 0.7064625585386706, 0.795197501037768,  0.8917384712110205, 0.8993440776696643, 43)
)

$rowIdx = 2
foreach ($entry in $data) {
    $ws.Cells.Item($rowIdx, 2).Value = $entry[1]
    $ws.Cells.Item($rowIdx, 3).Value = $entry[2]
    $ws.Cells.Item($rowIdx, 4).Value = $entry[3]
    $ws.Cells.Item($rowIdx, 5).Value = $entry[4]
    $ws.Cells.Item($rowIdx, 6).Value = $entry[5]
    $ws.Cells.Item($rowIdx, 7).Value = $entry[6]
    $rowIdx++
}

# Add new row 11 for Q9
$ws.Range("A11").Value = "Q9"
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Cells.Item(11, 2).Value = 0.08705267901537297
$ws.Cells.Item(11, 3).Value = 0.729473168299468
$ws.Cells.Item(11, 4).Value = 0.8438328393602245
$ws.Cells.Item(11, 5).Value = 0.9186037444732219
$ws.Cells.Item(11, 6).Value = 0.9255544964036263
$ws.Cells.Item(11, 7).Value = 42

$wb.Save()
